# Update player-stat cells on Sheet1.
# All touched cells are stored as text (t="inlineStr" in the source OOXML),
# holding digit strings like "2", "404", etc. Using a leading apostrophe
# keeps Excel's COM layer from reinterpreting the assigned text as a
# numeric value, so the cell keeps its original Text type.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "J2"  = "3"
    "E3"  = "494"
    "F3"  = "6"
    "G3"  = "5"
    "J7"  = "2"
    "E8"  = "1193"
    "F8"  = "16"
    "G8"  = "12"
    "E9"  = "441"
    "F9"  = "9"
    "H9"  = "5"
    "J9"  = "11"
    "E11" = "409"
    "F11" = "7"
    "G11" = "4"
    "L11" = "1"
    "E12" = "1419"
    "F12" = "16"
    "G12" = "16"
    "K12" = "2"
    "E13" = "672"
    "F13" = "9"
    "H13" = "2"
    "J13" = "2"
    "K13" = "1"
    "E14" = "361"
    "F14" = "6"
    "G14" = "4"
    "F16" = "16"
    "H16" = "7"
    "J16" = "7"
    "J17" = "12"
    "E18" = "870"
    "F18" = "14"
    "H18" = "3"
    "J18" = "4"
    "E22" = "116"
    "F22" = "2"
    "G22" = "2"
    "I22" = "2"
    "L22" = "1"
    "E23" = "458"
    "F23" = "11"
    "G23" = "5"
    "E25" = "1196"
    "F25" = "15"
    "G25" = "14"
    "I25" = "8"
    "E26" = "1243"
    "F26" = "16"
    "G26" = "14"
    "I26" = "7"
    "E27" = "964"
    "F27" = "16"
    "G27" = "12"
    "I27" = "9"
    "J28" = "4"
    "E29" = "551"
    "F29" = "14"
    "G29" = "5"
    "I29" = "5"
    "E30" = "18"
    "F30" = "2"
    "H30" = "2"
    "J30" = "3"
    "K30" = "1"
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = "'" + $updates[$addr]
}
